$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 72.5
$ws.Range("I9").Value = 45
$ws.Range("K9").Value = 45
$ws.Range("M9").Value = 124

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4855.8647
$ws.Range("I19").Value = 1555.1
$ws.Range("K19").Value = 1555.1
$ws.Range("M19").Value = -1380.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 186.28572
$ws.Range("J55").Value = 265.125
$ws.Range("L55").Value = 265.125
$ws.Range("N55").Value = -693.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2865
$ws.Range("I112").Value = 1718.8
$ws.Range("J112").Value = 3342.5833
$ws.Range("K112").Value = 5156.4
$ws.Range("L112").Value = 10027.7499
$ws.Range("M112").Value = -4048.4
$ws.Range("N112").Value = -12243.7499

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4038.5334
$ws.Range("I116").Value = 2899.75
$ws.Range("J116").Value = 5340
$ws.Range("K116").Value = 2899.75
$ws.Range("L116").Value = 5340
$ws.Range("M116").Value = 542.25
$ws.Range("N116").Value = -12224

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2599
$ws.Range("I135").Value = 2298.875
$ws.Range("J135").Value = 5000
$ws.Range("K135").Value = 20689.875
$ws.Range("L135").Value = 45000
$ws.Range("M135").Value = -18154.875
$ws.Range("N135").Value = -50070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3040.1904
$ws.Range("I63").Value = 2224.5
$ws.Range("J63").Value = 3781.7273
$ws.Range("K63").Value = 2224.5
$ws.Range("L63").Value = 3781.7273
$ws.Range("M63").Value = -1538.5
$ws.Range("N63").Value = -5153.7273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3040.1904
$ws.Range("I66").Value = 2224.5
$ws.Range("J66").Value = 3781.7273
$ws.Range("K66").Value = 11122.5
$ws.Range("L66").Value = 18908.6365
$ws.Range("M66").Value = -7690.5
$ws.Range("N66").Value = -25772.6365

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 49999
$ws.Range("J76").Value = 49999
$ws.Range("L76").Value = 49999
$ws.Range("N76").Value = -50675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 49999
$ws.Range("J79").Value = 49999
$ws.Range("L79").Value = 49999
$ws.Range("N79").Value = -52339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 99999
$ws.Range("J80").Value = 99999
$ws.Range("L80").Value = 99999
$ws.Range("N80").Value = -101995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 99999
$ws.Range("J83").Value = 99999
$ws.Range("L83").Value = 299997
$ws.Range("N83").Value = -309981

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1283
$ws.Range("J88").Value = 1371
$ws.Range("L88").Value = 1371
$ws.Range("N88").Value = -2183

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1283
$ws.Range("J91").Value = 1371
$ws.Range("L91").Value = 1371
$ws.Range("N91").Value = -4179

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 29442792
$ws.Range("I97").Value = 45455972
$ws.Range("K97").Value = 45455972
$ws.Range("M97").Value = -45455476

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1188.4242
$ws.Range("I94").Value = 838.25
$ws.Range("J94").Value = 1727.1538
$ws.Range("K94").Value = 838.25
$ws.Range("L94").Value = 1727.1538
$ws.Range("M94").Value = -387.25
$ws.Range("N94").Value = -2629.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2553.8865
$ws.Range("I134").Value = 2228.6099
$ws.Range("J134").Value = 6999.3335
$ws.Range("K134").Value = 6685.8297
$ws.Range("L134").Value = 20998.0005
$ws.Range("M134").Value = -4150.8297
$ws.Range("N134").Value = -26068.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 8867.857
$ws.Range("I22").Value = 23849.2
$ws.Range("J22").Value = 544.8889
$ws.Range("K22").Value = 23849.2
$ws.Range("L22").Value = 544.8889
$ws.Range("M22").Value = -23499.2
$ws.Range("N22").Value = -1244.8889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 2670.3333
$ws.Range("I25").Value = 2011
$ws.Range("K25").Value = 2011
$ws.Range("M25").Value = -1837

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2222
$ws.Range("I58").Value = 1771.6
$ws.Range("K58").Value = 1771.6
$ws.Range("M58").Value = -1568.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4796.25
$ws.Range("I122").Value = 1196.65
$ws.Range("J122").Value = 13795.25
$ws.Range("K122").Value = 3589.95
$ws.Range("L122").Value = 41385.75
$ws.Range("M122").Value = -1139.95
$ws.Range("N122").Value = -46285.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2222
$ws.Range("I136").Value = 1771.6
$ws.Range("K136").Value = 5314.799999999999
$ws.Range("M136").Value = -2764.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8821.053
$ws.Range("I3").Value = 6514.2856
$ws.Range("J3").Value = 10166.667
$ws.Range("K3").Value = 19542.8568
$ws.Range("L3").Value = 30500.001
$ws.Range("M3").Value = -19430.8568
$ws.Range("N3").Value = -30724.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3448.9412
$ws.Range("I139").Value = 2223.7144
$ws.Range("J139").Value = 9166.666999999999
$ws.Range("K139").Value = 6671.1432
$ws.Range("L139").Value = 27500.001
$ws.Range("M139").Value = -1531.1432
$ws.Range("N139").Value = -37780.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 19750
$ws.Range("J57").Value = 37500
$ws.Range("L57").Value = 37500
$ws.Range("N57").Value = -39140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 64476.832
$ws.Range("J136").Value = 64476.832
$ws.Range("L136").Value = 193430.496
$ws.Range("N136").Value = -198530.496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1639.2941
$ws.Range("I82").Value = 1849.7778
$ws.Range("J82").Value = 1402.5
$ws.Range("K82").Value = 1849.7778
$ws.Range("L82").Value = 1402.5
$ws.Range("M82").Value = -1488.7778
$ws.Range("N82").Value = -2124.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1639.2941
$ws.Range("I85").Value = 1849.7778
$ws.Range("J85").Value = 1402.5
$ws.Range("K85").Value = 1849.7778
$ws.Range("L85").Value = 1402.5
$ws.Range("M85").Value = -601.7778000000001
$ws.Range("N85").Value = -3898.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4200.2144
$ws.Range("I132").Value = 2970.8
$ws.Range("J132").Value = 10347.286
$ws.Range("K132").Value = 8912.400000000001
$ws.Range("L132").Value = 31041.858
$ws.Range("M132").Value = -6382.400000000001
$ws.Range("N132").Value = -36101.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3282.617
$ws.Range("I132").Value = 2855.6047
$ws.Range("J132").Value = 7873
$ws.Range("K132").Value = 8566.8141
$ws.Range("L132").Value = 23619
$ws.Range("M132").Value = -6036.8141
$ws.Range("N132").Value = -28679
